# Applies the "29/09/2017 CHITRA MAMTHA CHICK IN" purchase-details entry
# to the end of the existing ledger block in the document.
#
# Strategy: locate the last existing paragraph of the ledger
# ("Amount balance ... - 160951.0"), then append the eleven new
# paragraphs for the new visit record right after it, one at a time,
# using Range.InsertXML so each paragraph gets exact run/tab layout
# matching the style already used throughout the document (PlainText
# style, Courier New font, bold only on the blank spacer line and the
# final "Amount balance" line).

$d = $word.ActiveDocument

# ---- locate anchor paragraph -------------------------------------------
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*Amount balance*" -and $txt -like "*160951.0*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate the '- 160951.0' Amount balance paragraph"
}

# ---- XML building helpers ------------------------------------------------
$W_NS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-RPrXml([bool]$bold) {
    if ($bold) {
        return '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'
    } else {
        return '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
    }
}

function Escape-Xml([string]$s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    return $s
}

# Each run spec is a Hashtable @{ Tokens = @(...); Preserve = $true/$false }.
# Token "TAB" -> <w:tab/>; any other string -> <w:t>.
function New-RunSpec([array]$tokens, [bool]$preserve) {
    return @{ Tokens = $tokens; Preserve = $preserve }
}

function Build-RunXml([bool]$bold, [array]$tokens, [bool]$preserve) {
    $rpr = Get-RPrXml $bold
    $inner = ""
    foreach ($tok in $tokens) {
        if ($tok -eq "TAB") {
            $inner += "<w:tab/>"
        } else {
            $escaped = Escape-Xml $tok
            if ($preserve) {
                $inner += "<w:t xml:space=`"preserve`">$escaped</w:t>"
            } else {
                $inner += "<w:t>$escaped</w:t>"
            }
        }
    }
    return "<w:r>$rpr$inner</w:r>"
}

function Build-ParagraphXml([bool]$bold, [array]$runs) {
    $pPrRpr = Get-RPrXml $bold
    $pPr = "<w:pPr><w:pStyle w:val=`"PlainText`"/>$pPrRpr</w:pPr>"
    $runsXml = ""
    foreach ($run in $runs) {
        $runsXml += Build-RunXml $bold $run.Tokens $run.Preserve
    }
    return "<w:p>$pPr$runsXml</w:p>"
}

# Inserts a brand-new paragraph immediately after the paragraph at
# $afterIndex, with the given bold/run layout; returns the new
# paragraph's 1-based index.
function Insert-LedgerParagraph([int]$afterIndex, [bool]$bold, [array]$runs) {
    $anchor = $d.Paragraphs.Item($afterIndex)
    $anchor.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)

    $pXml = Build-ParagraphXml $bold $runs
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        "<pkg:xmlData><w:document $W_NS><w:body>$pXml</w:body></w:document></pkg:xmlData>" +
        '</pkg:part></pkg:package>'

    $newPara.Range.InsertXML($pkg)
    return $newIndex
}

# ---- the eleven new paragraphs -------------------------------------------
$cur = $anchorIndex

# 1) blank bold spacer line
$cur = Insert-LedgerParagraph $cur $true @()

# 2) timestamp line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("THU Sep 28") $false),
    (New-RunSpec @(" 13:20:55 PDT 2017") $true)
)

# 3) Person Name line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Person Name") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- M") $false)
)

# 4) separator dashed line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("---------------------------------------------------------------") $false)
)

# 5) Item Name line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Item Name") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- CARROT EVE") $false)
)

# 6) Number of Pockets line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Number of Pockets") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- 5") $false)
)

# 7) Number of KGs line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Number of KGs") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- 490") $false)
)

# 8) Rate line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Rate") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- 20") $false)
)

# 9) Transport & Miscellaneous line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Transport & Miscellaneous") $false),
    (New-RunSpec @("TAB", "- 50") $false)
)

# 10) Total Price line
$cur = Insert-LedgerParagraph $cur $false @(
    (New-RunSpec @("Total Price") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- 9850.0") $false)
)

# 11) Amount balance line (bold)
$cur = Insert-LedgerParagraph $cur $true @(
    (New-RunSpec @("Amount balance") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB") $false),
    (New-RunSpec @("TAB", "- 170801.0") $false)
)

Write-Host "Inserted ledger entry; final paragraph index:" $cur
Write-Host "Total paragraphs now:" $d.Paragraphs.Count
